$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2..19 (columns A..T), reflecting the updated TPM-based
# NATMI ligand-receptor computation. Columns A and D are cluster-name strings,
# B/C are constant ("Siglec1"/"Spn"), and E..T are the recomputed numeric
# statistics.
$rows = @(
    @("ECs", "Siglec1", "Spn", "Inflammatory-Mac", 1, 0.5, 0.1804155, 0.360831, 0.001449671343453619, 0.0009679009286790942, 3, 1, 2.019968666666667, 6.059906, 0.3745841340760022, 0.3745841340760022, 0.3644336569810001, 2.186601941886, 0.0005430238848823687, 0.0003625603312406168),
    @("ECs", "Siglec1", "Spn", "Neutrophils", 1, 0.5, 0.1804155, 0.360831, 0.001449671343453619, 0.0009679009286790942, 3, 1, 2.898691666666667, 8.696075, 0.537535025087018, 0.537535025087018, 0.5229689063875, 3.137813438325, 0.0007792491219712723, 0.0005202806499792649),
    @("ECs", "Siglec1", "Spn", "Resolving-Mac", 1, 0.5, 0.1804155, 0.360831, 0.001449671343453619, 0.0009679009286790942, 2, 0.6666666666666666, 0.473903, 1.421709, 0.08788084083697982, 0.08788084083697981, 0.0854994466965, 0.512996680179, 0.0001273983365999782, 0.00008505994745921242),
    @("FAPs", "Siglec1", "Spn", "Inflammatory-Mac", 1, 0.3333333333333333, 0.03714233333333333, 0.111427, 0.0002984454011009889, 0.0002988942102533469, 3, 1, 2.019968666666667, 6.059906, 0.3745841340760022, 0.3745841340760022, 0.07502634954022222, 0.675237145862, 0.0001117929121403791, 0.0001119610289280805),
    @("FAPs", "Siglec1", "Spn", "Neutrophils", 1, 0.3333333333333333, 0.03714233333333333, 0.111427, 0.0002984454011009889, 0.0002988942102533469, 3, 1, 2.898691666666667, 8.696075, 0.537535025087018, 0.537535025087018, 0.1076641721138889, 0.968977549025, 0.0001604248561679252, 0.0001606661068068973),
    @("FAPs", "Siglec1", "Spn", "Resolving-Mac", 1, 0.3333333333333333, 0.03714233333333333, 0.111427, 0.0002984454011009889, 0.0002988942102533469, 2, 0.6666666666666666, 0.473903, 1.421709, 0.08788084083697982, 0.08788084083697981, 0.01760186319366667, 0.158416768743, 0.00002622763279268461, 0.00002626707451836916),
    @("Inflammatory-Mac", "Siglec1", "Spn", "Inflammatory-Mac", 3, 1, 53.85241266666666, 161.557238, 0.4327139265678689, 0.4333646518592622, 3, 1, 2.019968666666667, 6.059906, 0.3745841340760022, 0.3745841340760022, 108.7801862110698, 979.0216758996279, 0.162087771486052, 0.1623315228558498),
    @("Inflammatory-Mac", "Siglec1", "Spn", "Neutrophils", 3, 1, 53.85241266666666, 161.557238, 0.4327139265678689, 0.4333646518592622, 3, 1, 2.898691666666667, 8.696075, 0.537535025087018, 0.537535025087018, 156.1015398267611, 1404.91385844085, 0.2325988913731615, 0.2329486790089953),
    @("Inflammatory-Mac", "Siglec1", "Spn", "Resolving-Mac", 3, 1, 53.85241266666666, 161.557238, 0.4327139265678689, 0.4333646518592622, 2, 0.6666666666666666, 0.473903, 1.421709, 0.08788084083697982, 0.08788084083697981, 25.52081991997133, 229.687379279742, 0.03802726370865546, 0.03808444999441698),
    @("MuSCs", "Siglec1", "Spn", "Inflammatory-Mac", 2, 1, 0.380206, 0.760412, 0.003055024334434274, 0.002039745700837033, 3, 1, 2.019968666666667, 6.059906, 0.3745841340760022, 0.3745841340760022, 0.7680042068786667, 4.608025241271999, 0.001144363644895177, 0.0007640563770832881),
    @("MuSCs", "Siglec1", "Spn", "Neutrophils", 2, 1, 0.380206, 0.760412, 0.003055024334434274, 0.002039745700837033, 3, 1, 2.898691666666667, 8.696075, 0.537535025087018, 0.537535025087018, 1.102099963816667, 6.6125997829, 0.001642182582251578, 0.001096434756470572),
    @("MuSCs", "Siglec1", "Spn", "Resolving-Mac", 2, 1, 0.380206, 0.760412, 0.003055024334434274, 0.002039745700837033, 2, 0.6666666666666666, 0.473903, 1.421709, 0.08788084083697982, 0.08788084083697981, 0.180180764018, 1.081084584108, 0.0002684781072875186, 0.0001792545672831731),
    @("Neutrophils", "Siglec1", "Spn", "Inflammatory-Mac", 3, 1, 14.52148466666667, 43.564454, 0.116682769416529, 0.1168582396856082, 3, 1, 2.019968666666667, 6.059906, 0.3745841340760022, 0.3745841340760022, 29.33294402014711, 263.996496181324, 0.04370751414348036, 0.04377324252227945),
    @("Neutrophils", "Siglec1", "Spn", "Neutrophils", 3, 1, 14.52148466666667, 43.564454, 0.116682769416529, 0.1168582396856082, 3, 1, 2.898691666666667, 8.696075, 0.537535025087018, 0.537535025087018, 42.09330659089444, 378.83975931805, 0.06272107538553667, 0.06281539680102814),
    @("Neutrophils", "Siglec1", "Spn", "Resolving-Mac", 3, 1, 14.52148466666667, 43.564454, 0.116682769416529, 0.1168582396856082, 2, 0.6666666666666666, 0.473903, 1.421709, 0.08788084083697982, 0.08788084083697981, 6.881775147987334, 61.935976331886, 0.01025417988751201, 0.01026960036230057),
    @("Resolving-Mac", "Siglec1", "Spn", "Inflammatory-Mac", 3, 1, 55.48103, 166.44309, 0.4458001629366132, 0.4464705676153602, 3, 1, 2.019968666666667, 6.059906, 0.3745841340760022, 0.3745841340760022, 112.0699421943933, 1008.62947974954, 0.166989668004552, 0.1672407909606209),
    @("Resolving-Mac", "Siglec1", "Spn", "Neutrophils", 3, 1, 55.48103, 166.44309, 0.4458001629366132, 0.4464705676153602, 3, 1, 2.898691666666667, 8.696075, 0.537535025087018, 0.537535025087018, 160.8223993190833, 1447.40159387175, 0.2396332017679291, 0.2399935677637378),
    @("Resolving-Mac", "Siglec1", "Spn", "Resolving-Mac", 3, 1, 55.48103, 166.44309, 0.4458001629366132, 0.4464705676153602, 2, 0.6666666666666666, 0.473903, 1.421709, 0.08788084083697982, 0.08788084083697981, 26.29262656009, 236.63363904081, 0.03917729316413218, 0.0392362088910015)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
